$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for added columns J and K
$ws.Range("J1").Value = "BranchAndReduce Duration [ns]"
$ws.Range("K1").Value = "BranchAndReduce Crossings"

# Updated timing columns (B, D, F, H) and new BranchAndReduce columns (J, K)
# Row 2 - matching_4_4.gr
$ws.Range("B2").Value = "17658"
$ws.Range("D2").Value = "2041"
$ws.Range("F2").Value = "2800"
$ws.Range("H2").Value = "2001"
$ws.Range("J2").Value = "14410"
$ws.Range("K2").Value = "0"

# Row 3 - cycle_8_sorted.gr
$ws.Range("B3").Value = "19801"
$ws.Range("D3").Value = "2888"
$ws.Range("F3").Value = "1119"
$ws.Range("H3").Value = "3094"
$ws.Range("J3").Value = "8152"
$ws.Range("K3").Value = "3"

# Row 4 - tree_6_10.gr
$ws.Range("B4").Value = "7286498594"
$ws.Range("D4").Value = "35228"
$ws.Range("F4").Value = "3201"
$ws.Range("H4").Value = "33148"
$ws.Range("J4").Value = "66427"
$ws.Range("K4").Value = "13"

# Row 5 - cycle_8_shuffled.gr
$ws.Range("B5").Value = "15386"
$ws.Range("D5").Value = "2794"
$ws.Range("F5").Value = "1099"
$ws.Range("H5").Value = "3043"
$ws.Range("J5").Value = "8603"
$ws.Range("K5").Value = "4"

# Row 6 - complete_4_5.gr
$ws.Range("B6").Value = "166373"
$ws.Range("D6").Value = "10707"
$ws.Range("F6").Value = "2920"
$ws.Range("H6").Value = "12638"
$ws.Range("J6").Value = "23942"
$ws.Range("K6").Value = "60"

# Row 7 - path_9_shuffled.gr
$ws.Range("B7").Value = "15270"
$ws.Range("D7").Value = "2708"
$ws.Range("F7").Value = "995"
$ws.Range("H7").Value = "11964"
$ws.Range("J7").Value = "11587"
$ws.Range("K7").Value = "6"

# Row 8 - ladder_4_4_sorted.gr
$ws.Range("B8").Value = "17763"
$ws.Range("D8").Value = "3393"
$ws.Range("F8").Value = "1246"
$ws.Range("H8").Value = "3780"
$ws.Range("J8").Value = "9698"
$ws.Range("K8").Value = "3"

# Row 9 - ladder_4_4_shuffled.gr
$ws.Range("B9").Value = "16871"
$ws.Range("D9").Value = "4863"
$ws.Range("F9").Value = "1246"
$ws.Range("H9").Value = "4822"
$ws.Range("J9").Value = "11855"
$ws.Range("K9").Value = "11"

# Row 10 - path_9_sorted.gr
$ws.Range("B10").Value = "13828"
$ws.Range("D10").Value = "2702"
$ws.Range("F10").Value = "1104"
$ws.Range("H10").Value = "2858"
$ws.Range("J10").Value = "10542"
$ws.Range("K10").Value = "0"

# Row 11 - website_20.gr
$ws.Range("B11").Value = "6106189208"
$ws.Range("D11").Value = "27576"
$ws.Range("F11").Value = "2919"
$ws.Range("H11").Value = "28088"
$ws.Range("J11").Value = "57681"
$ws.Range("K11").Value = "17"

# Row 12 - star_6.gr
$ws.Range("B12").Value = "440286"
$ws.Range("D12").Value = "5072"
$ws.Range("F12").Value = "1234"
$ws.Range("H12").Value = "5699"
$ws.Range("J12").Value = "14609"
$ws.Range("K12").Value = "0"

# Row 13 - plane_5_6.gr
$ws.Range("B13").Value = "624144"
$ws.Range("D13").Value = "7743"
$ws.Range("F13").Value = "1697"
$ws.Range("H13").Value = "7790"
$ws.Range("J13").Value = "20513"
$ws.Range("K13").Value = "0"

# Row 14 - grid_9_shuffled.gr
$ws.Range("B14").Value = "119308"
$ws.Range("D14").Value = "6271"
$ws.Range("F14").Value = "1664"
$ws.Range("H14").Value = "17429"
$ws.Range("J14").Value = "14656"
$ws.Range("K14").Value = "17"
